# add fixedtop settings for masthead
# Adds an "icon"/"slug" row (OVERVIEW) and a "slug" row (all topic sheets)
# to each sheet of the workbook, mirroring the existing key/value rows.

$wb = $excel.ActiveWorkbook

# OVERVIEW sheet: gets a new "icon" row (row 5) and a new "slug" row (row 6)
$ws = $wb.Worksheets.Item("OVERVIEW")
$ws.Range("A1:B1").Copy($ws.Range("A5:B5"))
$ws.Range("A5").Value = "icon"
$ws.Range("B5").Value = "fa-star"
$ws.Range("A1:B1").Copy($ws.Range("A6:B6"))
$ws.Range("A6").Value = "slug"
$ws.Range("B6").Value = "index"

# Each topic sheet already has rows 1-5 (topic, headline, byline,
# reporter_bio, icon) and now gets a new "slug" row appended as row 6.
$slugs = @{
    "HEALTH"     = "health-care"
    "PUBLIC_EDU" = "public-education"
    "HIGHER_EDU" = "higher-education"
    "TRANSPO"    = "transportation"
    "IMMIGRATION"= "immigration"
    "ENERGY"     = "energy"
    "ENVIRO"     = "environment"
    "TEF"        = "texas-enterprise-fund"
    "JUSTICE"    = "criminal-justice"
}

foreach ($name in $slugs.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1:B1").Copy($ws.Range("A6:B6"))
    $ws.Range("A6").Value = "slug"
    $ws.Range("B6").Value = $slugs[$name]
}
